$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Leetcode" row (row 4) problem counts
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 54

# Recalculate formulas (F4, C6, D6, F6 depend on these via SUM formulas)
$excel.Calculate()
